$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: fill in the previously-empty "Actual Start Date" (F15) ---
# Copy the date number format from the neighbouring E15 cell so the new
# cell picks up the same style index instead of inventing a new one, then
# write the value.
$ws.Range("E15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = 43892

# --- Row 17 (TASK013): follow-on schedule/status refresh bundled with
#     the TASK011 actual-start update above ---
# Assignation Date moves forward
$ws.Range("E17").Value = 43923
# Actual Start Date is cleared (task no longer shows an actual start yet)
$ws.Range("F17").ClearContents()
# Estimate/Dead line date pushed out
$ws.Range("G17").Value = 43953
# Actual Delivery Date is cleared as well
$ws.Range("I17").ClearContents()
# Status reverts from "done" to "Asseigned"
$ws.Range("J17").Value = "Asseigned"

# --- View state: move the selection to match the edited cell ---
$ws.Range("G17").Select()
